# Append two new data rows (23 and 24) to the NIFTY_Options_Analysis sheet,
# mirroring the formatting of the last existing data row (22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23 -----------------------------------------------------------
# Clone formatting (styles) from row 22 first, so the new rows pick up the
# same cell styles (AVOID red/maroon styles in C/D, numeric formats, etc.)
$ws.Range("A22:AE22").Copy()
$ws.Range("A23:AE23").PasteSpecial(-4122)

$ws.Range("A23").Value = "'2026-01-21"
$ws.Range("B23").Value = "'10:00:09"
$ws.Range("C23").Value = "AVOID"
$ws.Range("D23").Value = "AVOID"
$ws.Range("E23").Value = "'100%"
$ws.Range("F23").Value = "TRADEABLE"
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 25157.4
$ws.Range("I23").Value = 13.16
$ws.Range("J23").Value = 1.79
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 53.9
$ws.Range("M23").Value = "UNKNOWN"
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = "UNKNOWN"
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = 0
$ws.Range("S23").Value = 0
$ws.Range("T23").Value = "NONE"
$ws.Range("U23").Value = "'"
$ws.Range("V23").Value = 0
$ws.Range("W23").Value = 0
$ws.Range("X23").Value = 0
$ws.Range("Y23").Value = 0
$ws.Range("Z23").Value = 0
$ws.Range("AA23").Value = 0
$ws.Range("AB23").Value = 0
$ws.Range("AC23").Value = "HARD VETO: CPR TRENDING DAY: Price 25157.40 below BC 25378.17 - BEARISH TRENDING DAY likely"
$ws.Range("AD23").Value = "CPR TRENDING DAY: Price 25157.40 below BC 25378.17 - BEARISH TRENDING DAY likely"
$ws.Range("AE23").Value = "Yes"

# --- Row 24 -----------------------------------------------------------
$ws.Range("A22:AE22").Copy()
$ws.Range("A24:AE24").PasteSpecial(-4122)

$ws.Range("A24").Value = "'2026-01-22"
$ws.Range("B24").Value = "'10:00:08"
$ws.Range("C24").Value = "AVOID"
$ws.Range("D24").Value = "AVOID"
$ws.Range("E24").Value = "'100%"
$ws.Range("F24").Value = "TRADEABLE"
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 25397.4
$ws.Range("I24").Value = 13.47
$ws.Range("J24").Value = 1.64
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 56.7
$ws.Range("M24").Value = "UNKNOWN"
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = "UNKNOWN"
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = 0
$ws.Range("S24").Value = 0
$ws.Range("T24").Value = "NONE"
$ws.Range("U24").Value = "'"
$ws.Range("V24").Value = 0
$ws.Range("W24").Value = 0
$ws.Range("X24").Value = 0
$ws.Range("Y24").Value = 0
$ws.Range("Z24").Value = 0
$ws.Range("AA24").Value = 0
$ws.Range("AB24").Value = 0
$ws.Range("AC24").Value = "HARD VETO: CPR TRENDING DAY: Price 25397.40 above TC 25141.79 - BULLISH TRENDING DAY likely"
$ws.Range("AD24").Value = "CPR TRENDING DAY: Price 25397.40 above TC 25141.79 - BULLISH TRENDING DAY likely"
$ws.Range("AE24").Value = "Yes"
